# Add a new results row (player, 1) to the "Resultados" sheet,
# mirroring the existing rows in the table (e.g. A2:B27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "player"
$ws.Range("B28").Value = 1
